$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) / Volume(1h) (E) figures for the Sat Feb 4 2023 20:xx scrape refresh.
# Values are written with a leading "'" so Excel stores them as literal text (matching
# the workbook's existing inlineStr cells) instead of auto-converting to numbers/percentages.
$ws.Range("D2").Value = "'330.56"
$ws.Range("E2").Value = "'-0.16%"

$ws.Range("D3").Value = "'41.60"
$ws.Range("E3").Value = "'0.88%"

$ws.Range("D4").Value = "'5.700"
$ws.Range("E4").Value = "'0.40%"

$ws.Range("D5").Value = "'0.08390"
$ws.Range("E5").Value = "'4.09%"

$ws.Range("D6").Value = "'8.812"
$ws.Range("E6").Value = "'0.88%"

$ws.Range("D7").Value = "'1.993"
$ws.Range("E7").Value = "'-1.60%"

$ws.Range("D8").Value = "'4.476"
$ws.Range("E8").Value = "'-1.18%"

$ws.Range("E9").Value = "'-2.06%"

$ws.Range("D10").Value = "'0.9247"
$ws.Range("E10").Value = "'0.49%"

$ws.Range("D11").Value = "'0.1273"
$ws.Range("E11").Value = "'1.11%"

$ws.Range("D12").Value = "'0.1971"
$ws.Range("E12").Value = "'1.46%"

$ws.Range("D13").Value = "'0.09346"
$ws.Range("E13").Value = "'-0.39%"

$ws.Range("D14").Value = "'0.03885"
$ws.Range("E14").Value = "'4.74%"

$ws.Range("D15").Value = "'0.1061"
$ws.Range("E15").Value = "'0.48%"

$ws.Range("D16").Value = "'0.001309"
$ws.Range("E16").Value = "'0.82%"

$ws.Range("D17").Value = "'0.006106"
$ws.Range("E17").Value = "'-1.33%"

$ws.Range("D18").Value = "'3.425"
$ws.Range("E18").Value = "'1.78%"

$ws.Range("E19").Value = "'0.74%"

$ws.Range("D20").Value = "'8.925"
$ws.Range("E20").Value = "'7.66%"

$ws.Range("D21").Value = "'0.1362"
$ws.Range("E21").Value = "'-4.07%"

$ws.Range("D22").Value = "'0.2509"
$ws.Range("E22").Value = "'-5.61%"

$ws.Range("D23").Value = "'0.04398"
$ws.Range("E23").Value = "'-0.89%"

$ws.Range("D24").Value = "'0.001245"
$ws.Range("E24").Value = "'-1.37%"

$ws.Range("D25").Value = "'0.004394"
$ws.Range("E25").Value = "'2.22%"

$ws.Range("E26").Value = "'-3.98%"

$ws.Range("D27").Value = "'0.0003990"
$ws.Range("E27").Value = "'-0.07%"

$ws.Range("D39").Value = "'0.02821"
$ws.Range("E39").Value = "'-1.82%"

$ws.Range("D40").Value = "'0.05525"
$ws.Range("E40").Value = "'0.86%"

$ws.Range("E41").Value = "'2.35%"

$ws.Range("D42").Value = "'0.1438"
$ws.Range("E42").Value = "'1.45%"

$ws.Range("D43").Value = "'0.008974"
$ws.Range("E43").Value = "'-10.08%"

$ws.Range("E44").Value = "'-6.62%"

$ws.Range("E45").Value = "'-1.35%"

$ws.Range("D46").Value = "'0.00006927"
$ws.Range("E46").Value = "'2.10%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.05%"

$ws.Range("D48").Value = "'0.003461"
$ws.Range("E48").Value = "'14.86%"

$ws.Range("D49").Value = "'0.002278"
$ws.Range("E49").Value = "'-0.33%"

$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.05%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.05%"
